$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.045.17"
$ws.Range("E2").Value = "  -4.13%  "
$ws.Range("D3").Value = "1.961.06"
$ws.Range("E3").Value = "  -6.53%  "
$ws.Range("E4").Value = "  +0.69%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.54%  "
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4980"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4208"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.20"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08992"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.097"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.16%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.855"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.30%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.440"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.22%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.932.07"
$ws.Range("E15").Value = "  -4.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001096"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -9.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06666"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.007"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "  -6.80%  "
$ws.Range("D23").Value = "29.056.30"
$ws.Range("E23").Value = "  -4.04%  "
$ws.Range("E24").Value = "  -4.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.287"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.170"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -12.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.246"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -10.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.83%  "
$ws.Range("E31").Value = "  -8.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09826"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.529"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.788"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.695"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.49%  "
$ws.Range("E36").Value = "  -7.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.939"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -11.22%  "
$ws.Range("E38").Value = "  -4.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06285"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6420"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1980"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.007"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6174"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.167"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.279"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.465"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000323"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06874"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.84%  "
$ws.Range("E51").Value = "  -9.18%  "
